$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Xtras")

# The new "Xtras/OffloadQueue" sample rows are inserted right after the blank
# separator row (old row 4), pushing the existing EmailSubject/EmailTo/
# EmailTemplate/Orchestrator* rows down by 5 (matches the row-shift seen in
# the target sheet: old r5->r10, r6->r11, r7->r12, r9->r14, r10->r15,
# r11->r16, r12->r17). Inserting whole rows (rather than rewriting every
# cell) keeps the existing formatting/styles travelling with their content.
$ws.Range("A5:A9").EntireRow.Insert()

# The InputFolder row's description changes from the old AddDataToQueue
# text to the new OnloadQueue text.
$ws.Range("C3").Value = "XTRAS/OnloadQueue: Input data folder"

# Populate the freshly inserted rows with the new Offload/Onload queue
# sample config. Order matches the order new shared strings should be
# appended in (OutputPath row first, then OutputQueue row, then
# OutputTemplate row, then OutputSheet row).
$ws.Range("A6").Value = "OutputPath"
$ws.Range("B6").Value = "TestData/Output/[File]"
$ws.Range("C6").Value = "XTRAS/OffloadQueue: Output data file path"

$ws.Range("C5").Value = "XTRAS/OffloadQueue: Output queue to fetch data from"
$ws.Range("A5").Value = "OutputQueue"
$ws.Range("B5").Value = "RFW-ChorePile"

$ws.Range("A8").Value = "OutputTemplate"
$ws.Range("B8").Value = "Xtras/ExcelTemplate.xlsx"
$ws.Range("C8").Value = "XTRAS/OffloadQueue: Output Excel template file"

$ws.Range("A7").Value = "OutputSheet"
$ws.Range("B7").Value = "DataSheet"
$ws.Range("C7").Value = "XTRAS/OffloadQueue: Output Excel sheet name"

# Hyperlinks don't automatically follow their cell when rows are inserted
# above them in this engine, so drop the (now stale) links and re-add them
# at their new homes (B11 = EmailTo, B16 = OrchestratorUrl).
$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B11"), "mailto:love@uipath.com")
$ws.Range("B11").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B16"), "https://demo.uipath.com")
$ws.Range("B16").Style = "Hyperlink"

# Match the author's final selection.
$ws.Range("C8").Select()
